$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H18").Value = 310
$ws.Range("I18").Value = 296.66666
$ws.Range("K18").Value = 296.66666
$ws.Range("M18").Value = -12.66665999999998

$ws.Range("H88").Value = 4000
$ws.Range("I88").Value = 0
$ws.Range("J88").Value = 4000
$ws.Range("K88").Value = 0
$ws.Range("L88").Value = 4000
$ws.Range("N88").Value = -4812
$ws.Range("M88").ClearContents()

$ws.Range("H91").Value = 4000
$ws.Range("I91").Value = 0
$ws.Range("J91").Value = 4000
$ws.Range("K91").Value = 0
$ws.Range("L91").Value = 4000
$ws.Range("N91").Value = -6808
$ws.Range("M91").ClearContents()

$ws.Range("H98").Value = 1493.6428
$ws.Range("I98").Value = 1628
$ws.Range("K98").Value = 1628
$ws.Range("M98").Value = -130

$ws.Range("H122").Value = 1493.6428
$ws.Range("I122").Value = 1628
$ws.Range("K122").Value = 4884
$ws.Range("M122").Value = -2434

$ws.Range("H137").Value = 1613731.8
$ws.Range("I137").Value = 65361.92
$ws.Range("J137").Value = 2528677.5
$ws.Range("K137").Value = 196085.76
$ws.Range("L137").Value = 7586032.5
$ws.Range("M137").Value = -193535.76
$ws.Range("N137").Value = -7591132.5

$ws.Range("H138").Value = 4439.384
$ws.Range("I138").Value = 2999.75
$ws.Range("J138").Value = 4500
$ws.Range("K138").Value = 8999.25
$ws.Range("L138").Value = 13500
$ws.Range("M138").Value = -3859.25
$ws.Range("N138").Value = -23780

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 28782556
$ws.Range("I32").Value = 32976606
$ws.Range("K32").Value = 32976606
$ws.Range("M32").Value = -32976319

$ws.Range("H55").Value = 63063
$ws.Range("J55").Value = 69666
$ws.Range("L55").Value = 69666
$ws.Range("N55").Value = -70296

$ws.Range("H61").Value = 4794.2964
$ws.Range("I61").Value = 4702.1816
$ws.Range("J61").Value = 5199.6
$ws.Range("K61").Value = 4702.1816
$ws.Range("L61").Value = 5199.6
$ws.Range("M61").Value = -4490.1816
$ws.Range("N61").Value = -5623.6

$ws.Range("H74").Value = 3040.5
$ws.Range("J74").Value = 3250
$ws.Range("L74").Value = 3250
$ws.Range("N74").Value = -4998

$ws.Range("H77").Value = 3040.5
$ws.Range("J77").Value = 3250
$ws.Range("L77").Value = 16250
$ws.Range("N77").Value = -24986

$ws.Range("H88").Value = 2814.6667
$ws.Range("I88").Value = 2198
$ws.Range("J88").Value = 3123
$ws.Range("K88").Value = 2198
$ws.Range("L88").Value = 3123
$ws.Range("M88").Value = -1792
$ws.Range("N88").Value = -3935

$ws.Range("H91").Value = 2814.6667
$ws.Range("I91").Value = 2198
$ws.Range("J91").Value = 3123
$ws.Range("K91").Value = 2198
$ws.Range("L91").Value = 3123
$ws.Range("M91").Value = -794
$ws.Range("N91").Value = -5931

$ws.Range("H136").Value = 4794.2964
$ws.Range("I136").Value = 4702.1816
$ws.Range("J136").Value = 5199.6
$ws.Range("K136").Value = 14106.5448
$ws.Range("L136").Value = 15598.8
$ws.Range("M136").Value = -11556.5448
$ws.Range("N136").Value = -20698.8

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 4277.8936
$ws.Range("I31").Value = 2704.125
$ws.Range("J31").Value = 4600.718
$ws.Range("K31").Value = 2704.125
$ws.Range("L31").Value = 4600.718
$ws.Range("M31").Value = -2409.125
$ws.Range("N31").Value = -5190.718

$ws.Range("H34").Value = 4277.8936
$ws.Range("I34").Value = 2704.125
$ws.Range("J34").Value = 4600.718
$ws.Range("K34").Value = 2704.125
$ws.Range("L34").Value = 4600.718
$ws.Range("M34").Value = -2502.125
$ws.Range("N34").Value = -5004.718

$ws.Range("H58").Value = 3258.6584
$ws.Range("I58").Value = 2731.8235
$ws.Range("J58").Value = 3631.8333
$ws.Range("K58").Value = 2731.8235
$ws.Range("L58").Value = 3631.8333
$ws.Range("M58").Value = -2528.8235
$ws.Range("N58").Value = -4037.8333

$ws.Range("H59").Value = 0
$ws.Range("J59").Value = 0
$ws.Range("L59").Value = 0
$ws.Range("N59").ClearContents()

$ws.Range("H99").Value = 3081.4285
$ws.Range("I99").Value = 2845
$ws.Range("J99").Value = 3396.6667
$ws.Range("K99").Value = 2845
$ws.Range("L99").Value = 3396.6667
$ws.Range("M99").Value = -1347
$ws.Range("N99").Value = -6392.6667

$ws.Range("H126").Value = 3081.4285
$ws.Range("I126").Value = 2845
$ws.Range("J126").Value = 3396.6667
$ws.Range("K126").Value = 8535
$ws.Range("L126").Value = 10190.0001
$ws.Range("M126").Value = -6065
$ws.Range("N126").Value = -15130.0001

$ws.Range("H136").Value = 3258.6584
$ws.Range("I136").Value = 2731.8235
$ws.Range("J136").Value = 3631.8333
$ws.Range("K136").Value = 8195.470499999999
$ws.Range("L136").Value = 10895.4999
$ws.Range("M136").Value = -5645.470499999999
$ws.Range("N136").Value = -15995.4999

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H68").Value = 2414.093
$ws.Range("J68").Value = 2583.5
$ws.Range("L68").Value = 7750.5
$ws.Range("N68").Value = -9372.5

$ws.Range("H71").Value = 2414.093
$ws.Range("J71").Value = 2583.5
$ws.Range("L71").Value = 23251.5
$ws.Range("N71").Value = -31363.5

$ws.Range("H107").Value = 720.6429000000001
$ws.Range("I107").Value = 660.9231
$ws.Range("K107").Value = 1982.7693
$ws.Range("M107").Value = -62.76929999999993

$ws.Range("H122").Value = 775153.9
$ws.Range("I122").Value = 612.7143
$ws.Range("J122").Value = 3486048
$ws.Range("K122").Value = 5514.428699999999
$ws.Range("L122").Value = 31374432
$ws.Range("M122").Value = -3064.428699999999
$ws.Range("N122").Value = -31379332

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H126").Value = 2613.4
$ws.Range("I126").Value = 2133
$ws.Range("J126").Value = 3334
$ws.Range("K126").Value = 6399
$ws.Range("L126").Value = 10002
$ws.Range("M126").Value = -3929
$ws.Range("N126").Value = -14942

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H100").Value = 5634.3335
$ws.Range("I100").Value = 6001.5
$ws.Range("K100").Value = 6001.5
$ws.Range("M100").Value = -5460.5

$ws.Range("H129").Value = 183000
$ws.Range("J129").Value = 183000
$ws.Range("L129").Value = 183000
$ws.Range("N129").Value = -193000

$ws.Range("H136").Value = 5931.7666
$ws.Range("I136").Value = 4331.5557
$ws.Range("K136").Value = 12994.6671
$ws.Range("M136").Value = -10444.6671

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H81").Value = 3507.3333
$ws.Range("I81").Value = 2513.4
$ws.Range("J81").Value = 4749.75
$ws.Range("K81").Value = 5026.8
$ws.Range("L81").Value = 9499.5
$ws.Range("M81").Value = -3965.8
$ws.Range("N81").Value = -11621.5

$ws.Range("H84").Value = 3507.3333
$ws.Range("I84").Value = 2513.4
$ws.Range("J84").Value = 4749.75
$ws.Range("K84").Value = 25134
$ws.Range("L84").Value = 47497.5
$ws.Range("M84").Value = -19830
$ws.Range("N84").Value = -58105.5

$ws.Range("H136").Value = 52838.25
$ws.Range("I136").Value = 2653.8235
$ws.Range("J136").Value = 337216.66
$ws.Range("K136").Value = 7961.470499999999
$ws.Range("L136").Value = 1011649.98
$ws.Range("M136").Value = -5411.470499999999
$ws.Range("N136").Value = -1016749.98
